# Adds a new, blank slide at the end of the deck.
#
# This mirrors the author's edit: a new slide (final SlideID 269) was
# created at the end of the slide list. By the time the edit settled, the
# slide had no shapes on it (shapes that had briefly been added to it were
# removed again), so we simply append an empty slide using the built-in
# "Blank" layout.

$p = $ppt.ActivePresentation

# ppLayoutBlank = 12
$newIndex = $p.Slides.Count + 1
$newSlide = $p.Slides.Add($newIndex, 12)

Write-Output "Added slide $($newSlide.SlideIndex) with SlideID $($newSlide.SlideID); shapes=$($newSlide.Shapes.Count)"
